$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.103.19"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.302.10"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  -0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.69"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "97.43"
$ws.Range("E6").Value = "  -2.34%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.508"
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("E8").Value = "  -0.29%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.499"
$ws.Range("E9").Value = "  -1.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "33.65"
$ws.Range("E10").Value = "  -3.29%  "
$ws.Range("E11").Value = "  +0.51%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.26"
$ws.Range("E12").Value = "  -3.62%  "
$ws.Range("E13").Value = "  +2.22%  "
$ws.Range("E14").Value = "  +8.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.79"
$ws.Range("E15").Value = "  +0.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.661.03"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.315.52"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.809"
$ws.Range("E18").Value = "  +1.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "43.011.58"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0902"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.55"
$ws.Range("E21").Value = "  -1.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.05"
$ws.Range("E22").Value = "  +0.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.58"
$ws.Range("E23").Value = "  +0.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "236.68"
$ws.Range("E24").Value = "  +0.39%  "
$ws.Range("E25").Value = "  +2.76%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.45"
$ws.Range("E27").Value = "  -2.52%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.44"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.16"
$ws.Range("E29").Value = "  -0.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "166.73"
$ws.Range("E30").Value = "  +1.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "33.96"
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.15"
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.95"
$ws.Range("E34").Value = "  -1.10%  "
$ws.Range("B35").Value = "RenderToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.65"
$ws.Range("E35").Value = "  +4.24%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "16.94"
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("E37").Value = "  -1.43%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0693"
$ws.Range("E38").Value = "  -1.46%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.102"
$ws.Range("E39").Value = "  +0.81%  "
$ws.Range("B40").Value = "LidoDAOToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.82"
$ws.Range("E40").Value = "  -0.92%  "
$ws.Range("E41").Value = "  -2.70%  "
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("E43").Value = "  -3.67%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.985.68"
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("E45").Value = "  -0.69%  "
$ws.Range("E46").Value = "  +0.88%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "17.58"
$ws.Range("E47").Value = "  -4.08%  "
$ws.Range("E48").Value = "  -0.85%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.527.92"
$ws.Range("E49").Value = "  -0.56%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "53.21"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.60"
$ws.Range("E51").Value = "  -4.53%  "
